# Adjust pseudo code in Algorithm1:
# Each "for i <- 0 to n do" loop header (there are three of them, one per
# findThirdLargest pass) should read "for i <- 0 to n-1 do", matching the
# convention already used in Algorithm2 ("for I <- 0 to n - 1 do").
#
# Every Find/Replace below is scoped to a single paragraph's Range, and
# that paragraph is only touched when its text contains the exact literal
# "to n do" - which uniquely identifies the three Algorithm1 loop headers
# and leaves every other " n " occurrence in the document (e.g. "array A
# of n Integers", or the already-correct Algorithm2 header) untouched.

$d = $word.ActiveDocument

$count = $d.Paragraphs.Count
for ($i = 1; $i -le $count; $i++) {
    $para = $d.Paragraphs($i)
    $text = $para.Range.Text

    if ($text -like "*to n do*") {
        $rng = $para.Range
        $rng.Find.Execute(" n ", $true, $false, $false, $false, $false, $true, 1, $false, " n-1 ", 2)
    }
}
